$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure numeric-looking price strings stay as text (matches original inline-string cells)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.373.52"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "2.240.60"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "245.58"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "74.47"
$ws.Range("E7").Value = "  -3.98%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").Value = "43.55"
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("D11").Value = "0.0959"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "7.14"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "14.50"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "2.240.68"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "42.292.56"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "0.0000111"
$ws.Range("E18").Value = "  +12.63%  "
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").Value = "72.08"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "10.22"
$ws.Range("E21").Value = "  +40.03%  "
$ws.Range("D22").Value = "231.81"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -5.51%  "
$ws.Range("D24").Value = "11.80"
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  +4.09%  "
$ws.Range("D29").Value = "166.71"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").Value = "20.89"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").Value = "5.85"
$ws.Range("E31").Value = "  +18.95%  "
$ws.Range("D32").Value = "0.0814"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").Value = "30.01"
$ws.Range("E34").Value = "  -9.58%  "
$ws.Range("D35").Value = "0.124"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").Value = "0.0310"
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("D38").Value = "13.42"
$ws.Range("E38").Value = "  -6.63%  "
$ws.Range("D39").Value = "2.18"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "5.66"
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("D41").Value = "63.71"
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "8.86"
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("D44").Value = "105.92"
$ws.Range("E44").Value = "  -6.95%  "
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("D46").Value = "0.995"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("D51").Value = "4.12"
$ws.Range("E51").Value = "  -1.81%  "
